$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Insert a new "quality_flag" column (Q) right before the existing
# "notes" column (which slides from Q -> R).  We want the new column
# to pick up the same cell style (s="2", i.e. the style already used
# by the "alk_lab" column I) that the final XML shows, instead of the
# style that a plain insert-at-Q would inherit from its left neighbour
# (P, style "1").
#
# Trick: a whole-column Insert() copies the format of the column
# immediately to its LEFT.  So we first insert a throw-away blank
# column directly to the right of column I (which already uses style
# "2") to "clone" that exact style, then relocate that still-blank,
# still-styled column over to its real destination (immediately before
# "notes") using a column Cut+Insert (which preserves the style of the
# column being moved). This avoids creating any brand-new style entries
# in styles.xml - it just reuses the existing style index.
# ------------------------------------------------------------------

# Step 1: clone column I's style into a throwaway blank column at J.
# (this temporarily bumps everything from J rightwards by one column,
# including the "replicate_num" validation range J2:J11 -> K2:K11)
$ws.Columns.Item(10).Insert()

# Step 2: relocate that blank, style-cloned column from J (10) to its
# final resting place - immediately in front of "notes", which (after
# step 1's shift) now lives in column R (18).
$ws.Columns.Item(10).Cut()
$ws.Columns.Item(18).Insert()

# Step 3: undo the incidental side effect step 1 had on the
# "replicate_num" list validation (it got pushed from J2:J11 to
# K2:K11); put it back where it belongs.
$ws.Range("K2:K11").Validation.Delete()
$ws.Range("J2:J11").Validation.Add(3, 1, 1, '"1,2,3,4,5,6,7,8"')

# ------------------------------------------------------------------
# Populate the new quality_flag column (Q) - header + per-row values.
# ------------------------------------------------------------------
$ws.Range("Q1").Value = "quality_flag"

$ws.Range("Q2").Value = 2
$ws.Range("Q3").Value = 2
$ws.Range("Q4").Value = 2
$ws.Range("Q5").Value = 2
$ws.Range("Q6").Value = 2
$ws.Range("Q7").Value = 2
$ws.Range("Q8").Value = 2
$ws.Range("Q9").Value = 2
$ws.Range("Q10").Value = 1
$ws.Range("Q11").Value = 0

# Restore the sheet selection to where the author last left it.
$ws.Range("R14").Select()
